# Add two new submission rows to the first sheet
# ("八位序列号收集收集结果yd5"), which is the active sheet in this workbook.
#
# Row 135: Nov. / 2025-12-20 22:27:31 / 390e4e15 / 2113362931
# Row 136: 小天  / 2025-12-22 21:11:05 / 1480a607 / 603716468

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last populated row in column A so the new rows are appended
# right after the existing data (row 134).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$newRow1 = $lastRow + 1
$newRow2 = $lastRow + 2

$ws.Cells.Item($newRow1, 1).Value = "Nov."
$ws.Cells.Item($newRow1, 2).Value = 46011.935775463
$ws.Cells.Item($newRow1, 2).NumberFormat = "yyyy/m/d h:mm:ss;@"
$ws.Cells.Item($newRow1, 3).Value = "390e4e15"
$ws.Cells.Item($newRow1, 4).Value = "2113362931"

$ws.Cells.Item($newRow2, 1).Value = "小天"
$ws.Cells.Item($newRow2, 2).Value = 46013.8826967593
$ws.Cells.Item($newRow2, 2).NumberFormat = "yyyy/m/d h:mm:ss;@"
$ws.Cells.Item($newRow2, 3).Value = "1480a607"
$ws.Cells.Item($newRow2, 4).Value = "603716468"
